$d = $word.ActiveDocument

# Update the simulation run timestamps / duration shown in the report.

$d.Content.Find.Execute(
    "Start time: 2017-12-27 18:33:32", $true, $false, $false, $false, $false,
    $true, 1, $false, "Start time: 2018-01-31 12:37:07", 2
)

$d.Content.Find.Execute(
    "End time: 2017-12-27 18:33:51", $true, $false, $false, $false, $false,
    $true, 1, $false, "End time: 2018-01-31 12:37:29", 2
)

$d.Content.Find.Execute(
    "Duration: 19.73 secs", $true, $false, $false, $false, $false,
    $true, 1, $false, "Duration: 21.40 secs", 2
)
